$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("01-2015")

# --- Row 8: 采购订单 / 支付信息填写... -> 已解决, record date 2015-01-20 ---
$ws.Range("D8").Value = "已解决"
$ws.Range("E8").Value = 42024

# --- Row 9: 入库单 / 入库单明细里，最后多了一列 -> 已解决, record date ---
$ws.Range("D9").Value = "已解决"
$ws.Range("E9").Value = 42024

# --- Row 10: 入库单 / 入库单上的申请人，改为"操作人" -> 已解决, record date ---
$ws.Range("D10").Value = "已解决"
$ws.Range("E10").Value = 42024

# --- Row 11: 入库单 / 入库单，当货品入库... -> 已解决, record date ---
$ws.Range("D11").Value = "已解决"
$ws.Range("E11").Value = 42024

# --- Row 12: 大包申请单 / 部门选项... -> 无需修改, record date set, resolved date cleared ---
$ws.Range("D12").Value = "无需修改"
$ws.Range("E12").Value = 42024
$ws.Range("F12").ClearContents()

# --- Row 13: 大包协议维护 / 业务经理没有数据加载... -> 已解决, record date ---
$ws.Range("D13").Value = "已解决"
$ws.Range("E13").Value = 42024

# --- Row 14: 大包配送申请单 / 大包配送申请单，删除列... -> 无需修改, record date set, resolved date cleared ---
$ws.Range("D14").Value = "无需修改"
$ws.Range("E14").Value = 42024
$ws.Range("F14").ClearContents()

# --- Row 15: 大包配送申请单 / 为什么只能看到2种货品... -> stays 待解决, record date set, resolved date cleared ---
$ws.Range("E15").Value = 42024
$ws.Range("F15").ClearContents()

# --- Row 16: 采购订单 / 采购订单的货品没有按照帐套来筛选... -> text tightened, 已解决, record date ---
$ws.Range("C16").Value = "采购订单的货品没有按照帐套来筛选`n采购的货品列表应该根据帐套+供应商来筛选"
$ws.Range("D16").Value = "已解决"
$ws.Range("E16").Value = 42024

# --- Row 17: 大包申请单 / 具有审核权限的用户... -> text tightened, 无需修改, record date set, resolved date cleared ---
$ws.Range("C17").Value = "具有审核权限的用户，没有看到审核按钮，只有查看`n最后一列的空白列是什么？"
$ws.Range("D17").Value = "无需修改"
$ws.Range("E17").Value = 42024
$ws.Range("F17").ClearContents()

# --- View refresh: freeze header row, scroll to the newly updated rows, 80% zoom ---
$ws.Activate()
$ws.Range("C1").Select()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E16:F16").Select()
$excel.ActiveWindow.Zoom = 80
